$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 162, shifting existing rows 162:234 down to 163:235
$ws.Rows(162).Insert()

# Populate the newly inserted row 162 with the new record's data
$ws.Range("A162").Value = 9
$ws.Range("B162").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C162").Value = "Metropolitana"
$ws.Range("D162").Value = 44704
$ws.Range("E162").Value = 13
$ws.Range("F162").Value = 100112003
$ws.Range("G162").Value = "Ajo"
$ws.Range("H162").Value = "Chino"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 350
$ws.Range("K162").Value = 17000
$ws.Range("L162").Value = 18000
$ws.Range("M162").Value = 17571
$ws.Range("N162").Value = "`$/caja 10 kilos"
$ws.Range("O162").Value = "China"
$ws.Range("P162").Value = 1757
$ws.Range("Q162").Value = 10
$ws.Range("R162").Value = "Hortaliza"
